# Update Name of Algo
# Apply updated RandomForest imputation results to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -12.0843
$ws.Range("B3").Value  = 5.980799999999989
$ws.Range("C5").Value  = -14.47340000000001
$ws.Range("D5").Value  = -8.486699999999997
$ws.Range("D9").Value  = -8.683300000000004
$ws.Range("D11").Value = -8.3429
$ws.Range("B14").Value = 8.9322
$ws.Range("B21").Value = 5.661099999999993
$ws.Range("D21").Value = -7.671600000000003
$ws.Range("B23").Value = 5.922399999999999
$ws.Range("B25").Value = 5.813399999999993
